{"js": "// Supplier agreement text updates (Juno / OLPRODLOC).\n// For each (find, replace) pair, search the body for the exact\n// existing text and replace it in place, preserving the run's\n// formatting (insertText with InsertLocation.replace keeps the\n// surrounding run/paragraph intact).\nconst replacements = [\n  [\n    \"\u30ce\u30fc\u30b9\u30a6\u30a3\u30f3\u30c9\u30c8\u30ec\u30fc\u30c0\u30fc\u30ba\u306f\u3001\u30ef\u30a4\u30c9\u30ef\u30fc\u30eb\u30c9\u30a4\u30f3\u30dd\u30fc\u30bf\u30fc\u306e\u305f\u3081\u306e\u30d3\u30fc\u30eb\u3068\u91b8\u9020\u6240\u306e\u512a\u5148\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u3067\u3059\u3002\",\n    \"Northwind Traders \u306f\u3001Wide World Importers \u306e\u305f\u3081\u306e\u30d3\u30fc\u30eb\u3068\u30b5\u30a4\u30c0\u30fc\u306e\u512a\u5148\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u3067\u3059\u3002\",\n  ],\n  [\n    \"2023 \u5e74 2 \u6708 1 \u65e5\u306b Wide World Importers \u3068\u4ea4\u6e09\u3055\u308c\u305f\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u5951\u7d04\u6761\u4ef6\u306f\u6b21\u306e\u3068\u304a\u308a\u3067\u3059\u3002\",\n    \"2023 \u5e74 2 \u6708 1 \u65e5\u306b Wide World Importers \u3068\u4ea4\u6e09\u3057\u305f\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u5951\u7d04\u6761\u4ef6\u306f\u6b21\u306e\u3068\u304a\u308a\u3067\u3059\u3002\",\n  ],\n  [\n    \"10 \u65e5\u4ee5\u5185\u652f\u6255\u5272\u5f15 2%\u3001\u652f\u6255\u671f\u9650 30 \u65e5\u4ee5\u5185\",\n    \"10 \u65e5\u4ee5\u5185\u652f\u6255\u5272\u5f15 2%\u3001\u652f\u6255\u671f\u9650 45 \u65e5\u4ee5\u5185\",\n  ],\n  [\"\u4e00\u5f8b\u6599\u91d1\", \"\u9045\u5ef6\u6599\u91d1\"],\n  [\"1 \u304b\u6708\u3042\u305f\u308a $100\", \"1 \u304b\u6708\u3042\u305f\u308a 2%\"],\n  [\"\u6700\u5c0f\u6ce8\u6587\u91d1\u984d\", \"\u6700\u5c0f\u6ce8\u6587\u91cf\"],\n  [\"20 \u6642\u9593/\u6708\", \"1 \u304b\u6708\u3042\u305f\u308a 50 \u30b1\u30fc\u30b9\"],\n  [\"\u6700\u5927\u6ce8\u6587\u91d1\u984d\", \"\u6700\u5927\u6ce8\u6587\u91cf\"],\n  [\"0 (\u6700\u5927\u5024\u306a\u3057)\", \"\u6700\u5927\u5024\u306a\u3057\"],\n  [\n    \"\u30b1\u30fc\u30b9\u3042\u305f\u308a\u306e\u4fa1\u683c\u306f $25 \u3067\u56fa\u5b9a\u3055\u308c\u3066\u3044\u307e\u3059\",\n    \"\u30b1\u30fc\u30b9\u3042\u305f\u308a\u306e\u4fa1\u683c\u306f 25 \u30c9\u30eb\u3067\u56fa\u5b9a\u3055\u308c\u3066\u3044\u307e\u3059\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Supplier agreement text updates (Juno / OLPRODLOC).\n# Each pair is an exact existing string and its replacement; every\n# target string is unique in the document, so a plain Find/Replace is\n# unambiguous.\n$replacements = @(\n    @(\"\u30ce\u30fc\u30b9\u30a6\u30a3\u30f3\u30c9\u30c8\u30ec\u30fc\u30c0\u30fc\u30ba\u306f\u3001\u30ef\u30a4\u30c9\u30ef\u30fc\u30eb\u30c9\u30a4\u30f3\u30dd\u30fc\u30bf\u30fc\u306e\u305f\u3081\u306e\u30d3\u30fc\u30eb\u3068\u91b8\u9020\u6240\u306e\u512a\u5148\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u3067\u3059\u3002\", \"Northwind Traders \u306f\u3001Wide World Importers \u306e\u305f\u3081\u306e\u30d3\u30fc\u30eb\u3068\u30b5\u30a4\u30c0\u30fc\u306e\u512a\u5148\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u3067\u3059\u3002\"),\n    @(\"2023 \u5e74 2 \u6708 1 \u65e5\u306b Wide World Importers \u3068\u4ea4\u6e09\u3055\u308c\u305f\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u5951\u7d04\u6761\u4ef6\u306f\u6b21\u306e\u3068\u304a\u308a\u3067\u3059\u3002\", \"2023 \u5e74 2 \u6708 1 \u65e5\u306b Wide World Importers \u3068\u4ea4\u6e09\u3057\u305f\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u5951\u7d04\u6761\u4ef6\u306f\u6b21\u306e\u3068\u304a\u308a\u3067\u3059\u3002\"),\n    @(\"10 \u65e5\u4ee5\u5185\u652f\u6255\u5272\u5f15 2%\u3001\u652f\u6255\u671f\u9650 30 \u65e5\u4ee5\u5185\", \"10 \u65e5\u4ee5\u5185\u652f\u6255\u5272\u5f15 2%\u3001\u652f\u6255\u671f\u9650 45 \u65e5\u4ee5\u5185\"),\n    @(\"\u4e00\u5f8b\u6599\u91d1\", \"\u9045\u5ef6\u6599\u91d1\"),\n    @(\"1 \u304b\u6708\u3042\u305f\u308a `$100\", \"1 \u304b\u6708\u3042\u305f\u308a 2%\"),\n    @(\"\u6700\u5c0f\u6ce8\u6587\u91d1\u984d\", \"\u6700\u5c0f\u6ce8\u6587\u91cf\"),\n    @(\"20 \u6642\u9593/\u6708\", \"1 \u304b\u6708\u3042\u305f\u308a 50 \u30b1\u30fc\u30b9\"),\n    @(\"\u6700\u5927\u6ce8\u6587\u91d1\u984d\", \"\u6700\u5927\u6ce8\u6587\u91cf\"),\n    @(\"0 (\u6700\u5927\u5024\u306a\u3057)\", \"\u6700\u5927\u5024\u306a\u3057\"),\n    @(\"\u30b1\u30fc\u30b9\u3042\u305f\u308a\u306e\u4fa1\u683c\u306f `$25 \u3067\u56fa\u5b9a\u3055\u308c\u3066\u3044\u307e\u3059\", \"\u30b1\u30fc\u30b9\u3042\u305f\u308a\u306e\u4fa1\u683c\u306f 25 \u30c9\u30eb\u3067\u56fa\u5b9a\u3055\u308c\u3066\u3044\u307e\u3059\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.Execute([ref]$find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\n# The opening paragraph holds two sentences as two separate runs that\n# happen to share identical run formatting. Editing either sentence's\n# text through Find/Replace re-flows the paragraph and folds the two\n# runs into one (Word's own editor does the same coalescing when\n# formatting is identical). Nudge the second sentence's font color\n# away and back to its original value (wdColorBlack, matching the\n# existing <w:color w:val=\"000000\"/>) so Word re-splits the paragraph\n# back into two runs at the sentence boundary, matching the original\n# paragraph structure.\n$introPara = $d.Paragraphs.Item(3)\n$secondSentence = \"2023 \u5e74 2 \u6708 1 \u65e5\u306b Wide World Importers \u3068\u4ea4\u6e09\u3057\u305f\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u5951\u7d04\u6761\u4ef6\u306f\u6b21\u306e\u3068\u304a\u308a\u3067\u3059\u3002\"\n$splitStart = $introPara.Range.Start + (\"Northwind Traders \u306f\u3001Wide World Importers \u306e\u305f\u3081\u306e\u30d3\u30fc\u30eb\u3068\u30b5\u30a4\u30c0\u30fc\u306e\u512a\u5148\u30b5\u30d7\u30e9\u30a4\u30e4\u30fc\u3067\u3059\u3002\").Length\n$splitEnd = $introPara.Range.End - 1\n$secondRange = $d.Range($splitStart, $splitEnd)\nif ($secondRange.Text -eq $secondSentence) {\n    $secondRange.Font.Color = 1\n    $secondRange.Font.Color = 0\n}\n"}
